$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rename Sheet1 -> "Portfolio 1" ---
$ws1.Name = "Portfolio 1"

# --- Create "Portfolio 2" and "Port3" as copies of Portfolio 1 so they ---
# --- inherit the same row-height / formatting metadata, then overwrite ---
# --- their contents.                                                    ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Portfolio 2"

$ws1.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Port3"

# --- Portfolio 2 data ---
$ws2.Range("A2").Value = "AMZN"
$ws2.Range("B2").Value = 0.05
$ws2.Range("A3").Value = "FB"
$ws2.Range("B3").Value = 0.1
$ws2.Range("A4").Value = "AMZN"
$ws2.Range("B4").Value = 0.1
$ws2.Range("A5").Value = "AAPL"
$ws2.Range("B5").Value = 0.1
$ws2.Range("A6").Value = "MSFT"
$ws2.Range("B6").Value = 0.05
$ws2.Range("A7").Value = "GOOGL"
$ws2.Range("B7").Value = 0.1
$ws2.Range("A8").Value = "AMD"
$ws2.Range("B8").Value = 0.1
$ws2.Range("A9").Value = "CRM"
$ws2.Range("B9").Value = 0.1
$ws2.Range("A10").Value = "NVDA"
$ws2.Range("B10").Value = 0.1
$ws2.Range("A11").Value = "LRCX"
$ws2.Range("B11").Value = 0.1
$ws2.Range("A12").Value = "PYPL"
$ws2.Range("B12").Value = 0.1
$ws2.Range("B2:B12").NumberFormat = "0%"

# --- Port3 data (ticker list only, column A) ---
$ws3.Range("A2").Value = "AMZN"
$ws3.Range("A3").Value = "FB"
$ws3.Range("A4").Value = "AMZN"
$ws3.Range("A5").Value = "AAPL"
$ws3.Range("A6").Value = "LOWE"
$ws3.Range("A7").Value = "SHOP"
$ws3.Range("A8").Value = "CMG"
$ws3.Range("A9").Value = "ETSY"
$ws3.Range("A10").Value = "COST"
$ws3.Range("A11").Value = "QCOM"
$ws3.Range("A12").Value = "YETI"
$ws3.Range("B1:B7").Clear()

# --- Sheet view / selection updates ---
$ws1.Range("A1:B7").Select()
$ws2.Range("A2:A12").Select()
$ws3.Range("D4").Select()

Write-Host "done"
